# Update cryptocurrency price/volume data per upstream scrape refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "'25.853.42"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "'1.741.47"
$ws.Range("E3").Value = "  -0.91%  "
$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'225.41"
$ws.Range("E5").Value = "  -5.07%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'0.5144"
$ws.Range("E7").Value = "  +1.36%  "
$ws.Range("D8").Value = "'0.2789"
$ws.Range("E8").Value = "  +5.08%  "
$ws.Range("D9").Value = "'38.97"
$ws.Range("E9").Value = "  -5.42%  "
$ws.Range("D10").Value = "'0.06092"
$ws.Range("E10").Value = "  -1.90%  "
$ws.Range("D11").Value = "'1.740.85"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").Value = "'0.06979"
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'15.21"
$ws.Range("E13").Value = "  -2.35%  "
$ws.Range("D14").Value = "'0.6331"
$ws.Range("E14").Value = "  +4.45%  "
$ws.Range("D15").Value = "'4.496"
$ws.Range("E15").Value = "  +0.85%  "
$ws.Range("D16").Value = "'76.50"
$ws.Range("E16").Value = "  -1.38%  "
$ws.Range("D17").Value = "'1.000"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "'0.9999"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("D19").Value = "'25.879.98"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  -2.06%  "
$ws.Range("D21").Value = "'0.000006580"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("D22").Value = "'1.959.39"
$ws.Range("E22").Value = "  -0.91%  "
$ws.Range("D23").Value = "'4.074"
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'8.445"
$ws.Range("E24").Value = "  +3.68%  "
$ws.Range("D25").Value = "'5.093"
$ws.Range("E25").Value = "  -1.76%  "
$ws.Range("D26").Value = "'137.84"
$ws.Range("E26").Value = "  -0.01%  "
$ws.Range("D27").Value = "'1.509"
$ws.Range("E27").Value = "  +3.46%  "
$ws.Range("D28").Value = "'1.809"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("D29").Value = "'14.96"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "'102.61"
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").Value = "'0.08266"
$ws.Range("E31").Value = "  +0.36%  "
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("D33").Value = "'3.399"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("D34").Value = "'0.04405"
$ws.Range("E34").Value = "  +0.64%  "
$ws.Range("D35").Value = "'2.618"
$ws.Range("E35").Value = "  -1.37%  "
$ws.Range("D36").Value = "'0.9706"
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("D37").Value = "'0.5977"
$ws.Range("E37").Value = "  -0.53%  "
$ws.Range("D38").Value = "'2.670"
$ws.Range("E38").Value = "  -2.33%  "
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").Value = "'0.9990"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("E41").Value = "  -1.74%  "
$ws.Range("D42").Value = "'100.58"
$ws.Range("E42").Value = "  -2.60%  "
$ws.Range("D43").Value = "'0.3816"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("D44").Value = "'0.7267"
$ws.Range("E44").Value = "  -2.43%  "
$ws.Range("D45").Value = "'4.870"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").Value = "'0.05463"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").Value = "'6.253"
$ws.Range("E47").Value = "  +4.76%  "
$ws.Range("D48").Value = "'0.1102"
$ws.Range("E48").Value = "  +1.84%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'29.64"
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").Value = "'52.05"
$ws.Range("E50").Value = "  -0.09%  "
$ws.Range("D51").Value = "'7.467"
$ws.Range("E51").Value = "  -0.53%  "
